# Apply "new rail car" trace report update:
# - Add a 4th trace-event data row (CGAX10167 moves from row 3 to row 6)
# - Update/replace the other three rows with new trip data (new times, weights, etc.)
# - Update the summary line in A1 (new completion date/time + event count)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: search summary text ---
$ws.Range("A1").Value = "Description unknown, completed 06/21/2023 08:46:33 EDT, by WPJTOWN1.The search returned: 4 events."

# --- Row 3: BNGX 30727 ---
$ws.Range("A3").Value = "BNGX"
$ws.Range("B3").Value = 30727
$ws.Range("C3").Value = "DENVER"
$ws.Range("D3").Value = "CO"
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 19
$ws.Range("G3").Value = 1259
$ws.Range("H3").Value = "Arrive In-Transit"
$ws.Range("I3").Value = "HKCKDE"
$ws.Range("J3").Value = "LOVELAND"
$ws.Range("K3").Value = "CO"
$ws.Range("L3").Value = 282200
$ws.Range("M3").Value = 64400
$ws.Range("N3").Value = 217800
$ws.Range("O3").Value = "BNGX30727"

# --- Row 4: CAIX 541012 ---
$ws.Range("A4").Value = "CAIX"
$ws.Range("B4").Value = 541012
$ws.Range("C4").Value = "HUTCHINSON"
$ws.Range("D4").Value = "KS"
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 21
$ws.Range("G4").Value = 357
$ws.Range("H4").Value = "Departure"
$ws.Range("I4").Value = "HKCKDE"
$ws.Range("J4").Value = "LOVELAND"
$ws.Range("K4").Value = "CO"
$ws.Range("L4").Value = 273100
$ws.Range("M4").Value = 62900
$ws.Range("N4").Value = 210200
$ws.Range("O4").Value = "CAIX541012"

# --- Row 5: CGEX 1941 ---
$ws.Range("A5").Value = "CGEX"
$ws.Range("B5").Value = 1941
$ws.Range("C5").Value = "LA SALLE"
$ws.Range("D5").Value = "CO"
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 20
$ws.Range("G5").Value = 1855
$ws.Range("H5").Value = "Arrive In-Transit"
$ws.Range("I5").Value = "LDI602"
$ws.Range("J5").Value = "JOHNSTOWN"
$ws.Range("K5").Value = "CO"
$ws.Range("L5").Value = 198750
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 198750
$ws.Range("O5").Value = "CGEX1941"

# --- Row 6 (new): CGAX 10167 ---
$ws.Range("A6").Value = "CGAX"
$ws.Range("B6").Value = 10167
$ws.Range("C6").Value = "WINDSOR"
$ws.Range("D6").Value = "CO"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 16
$ws.Range("G6").Value = 1524
$ws.Range("H6").Value = "Arrive In-Transit"
$ws.Range("J6").Value = "JOHNSTOWN"
$ws.Range("K6").Value = "CO"
$ws.Range("L6").Value = 273000
$ws.Range("M6").Value = 64200
$ws.Range("N6").Value = 208800
$ws.Range("O6").Value = "CGAX10167"

# Update the selection to mirror the new data extent
$ws.Range("O3:O6").Select()
